# Apply the edits described by the commit "compare the metrics and plot distribution"
# - Move the active selection from C42 to B29 (and scroll the view toward A43)
# - Give rows 17 and 18 an explicit 16.5pt custom row height (row 18 was
#   previously squashed to 0.4pt; row 17 gets the same new custom height)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the visible/scrolled area and the active cell selection on the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("B29").Select()

# Rows 17 and 18 now both use a custom row height of 16.5 points.
$ws.Rows.Item(17).RowHeight = 16.5
$ws.Rows.Item(18).RowHeight = 16.5
